$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values per diff (future production now uses uncon_planned_qty)
$ws.Range("G2").Value = 660
$ws.Range("H2").Value = 660
$ws.Range("J2").Value = 621

# Add new row 3 for MAT_B / LINE_B
$ws.Range("A3").Value = "MAT_B"
$ws.Range("B3").Value = "PLANT_001"
$ws.Range("C3").Value = "LINE_B"

$ws.Range("D3").Value = 45294
$ws.Range("E3").Value = 45295
$ws.Range("F3").Value = 45296
$ws.Range("D3:F3").Style = $ws.Range("D2:F2").Style
$ws.Range("D3:F3").NumberFormat = $ws.Range("D2:F2").NumberFormat

$ws.Range("G3").Value = 80
$ws.Range("H3").Value = 80
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = 75
